$wb = $excel.ActiveWorkbook

$sheetNames = @("Add Devices Loop A", "Other Devices Loop A")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Copy formatting for the new "Loop"/"Column" header cells from the
    # existing "DC Unit Loading Details Name" header cell (I1), which already
    # carries the bold/filled/bordered header style.
    $ws.Range("I1").Copy()
    $ws.Range("J1:K1").PasteSpecial(-4122)

    # --- Copy formatting for the new "Built-in Loop-*" cells from B4, which
    # currently carries the plain bordered style that these new cells adopt.
    $ws.Range("B4").Copy()
    $ws.Range("J2:J5").PasteSpecial(-4122)
    $ws.Range("K2").PasteSpecial(-4122)

    $excel.CutCopyMode = 0

    # --- New header values.
    $ws.Range("J1").Value = "Loop"
    $ws.Range("K1").Value = "Column"

    # --- New built-in loop rows.
    $ws.Range("J2").Value = "Built-in Loop-A"
    $ws.Range("K2").Value = 2
    $ws.Range("J3").Value = "Built-in Loop-B"
    $ws.Range("J4").Value = "Built-in Loop-C"
    $ws.Range("J5").Value = "Built-in Loop-D"

    # --- B4 becomes a plain (unformatted) cell holding the user-story note.
    $ws.Range("B4").ClearFormats()
    $ws.Range("B4").Value = "NGC-1826/T918 OR TC-63797"
}

# --- Updated DC Units values (340.8 -> 341).
$wsAdd = $wb.Worksheets.Item("Add Devices Loop A")
$wsAdd.Range("G1").Value = 341

$wsOther = $wb.Worksheets.Item("Other Devices Loop A")
$wsOther.Range("G2").Value = 341

$wsAdd.Range("B4").Select()
$wsOther.Range("B4").Select()
